$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 31629
$ws.Range("B3").Value = 92.06928040722123
$ws.Range("B4").Value = 10.1529448215195
$ws.Range("B5").Value = 47.44
$ws.Range("B6").Value = 85.31
$ws.Range("B7").Value = 97.39
